$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their original text representation
# (values like "0.400", "165.60" must not be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '65.335.69'
$ws.Range('E2').Value = '  +3.52%  '
$ws.Range('D3').Value = '2.644.25'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('D5').Value = '603.57'
$ws.Range('E5').Value = '  +2.69%  '
$ws.Range('D6').Value = '156.35'
$ws.Range('E6').Value = '  +4.49%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '0.594'
$ws.Range('E8').Value = '  +1.81%  '
$ws.Range('E9').Value = '  +8.33%  '
$ws.Range('D10').Value = '0.400'
$ws.Range('E10').Value = '  +4.49%  '
$ws.Range('D11').Value = '5.82'
$ws.Range('E11').Value = '  +3.25%  '
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').Value = '29.02'
$ws.Range('E13').Value = '  +6.19%  '
$ws.Range('D14').Value = '0.0000186'
$ws.Range('E14').Value = '  +20.03%  '
$ws.Range('D15').Value = '3.118.28'
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('D16').Value = '65.233.00'
$ws.Range('E16').Value = '  +3.54%  '
$ws.Range('D17').Value = '2.784.24'
$ws.Range('E17').Value = '  +8.05%  '
$ws.Range('D18').Value = '12.59'
$ws.Range('E18').Value = '  +2.57%  '
$ws.Range('D19').Value = '4.81'
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').Value = '355.18'
$ws.Range('E20').Value = '  +2.98%  '
$ws.Range('D21').Value = '7.26'
$ws.Range('E21').Value = '  +5.79%  '
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').Value = '68.18'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('D24').Value = '1.72'
$ws.Range('E24').Value = '  +2.06%  '
$ws.Range('D25').Value = '9.52'
$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('D26').Value = '1.66'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = '8.14'
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('D28').Value = '0.164'
$ws.Range('E28').Value = '  +2.49%  '
$ws.Range('D29').Value = '0.0₃0957'
$ws.Range('E29').Value = '  +13.55%  '
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('D31').Value = '2.12'
$ws.Range('E31').Value = '  +4.15%  '
$ws.Range('D32').Value = '522.18'
$ws.Range('E32').Value = '  -4.66%  '
$ws.Range('D33').Value = '1.79'
$ws.Range('E33').Value = '  +2.98%  '
$ws.Range('D34').Value = '5.67'
$ws.Range('E34').Value = '  +9.68%  '
$ws.Range('D35').Value = '6.38'
$ws.Range('E35').Value = '  +4.80%  '
$ws.Range('D36').Value = '0.428'
$ws.Range('E36').Value = '  +4.18%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '165.60'
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '2.04'
$ws.Range('E38').Value = '  +6.71%  '
$ws.Range('D39').Value = '20.22'
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = '42.13'
$ws.Range('E42').Value = '  +6.53%  '
$ws.Range('D43').Value = '165.75'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').Value = '4.10'
$ws.Range('E44').Value = '  +4.42%  '
$ws.Range('D45').Value = '0.0609'
$ws.Range('E45').Value = '  +4.88%  '
$ws.Range('D46').Value = '23.20'
$ws.Range('E46').Value = '  +3.61%  '
$ws.Range('D47').Value = '2.21'
$ws.Range('E47').Value = '  +8.81%  '
$ws.Range('D48').Value = '0.653'
$ws.Range('E48').Value = '  +4.01%  '
$ws.Range('D49').Value = '0.0254'
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('D50').Value = '0.0986'
$ws.Range('E50').Value = '  +2.72%  '
$ws.Range('D51').Value = '19.55'
$ws.Range('E51').Value = '  +2.93%  '
